# Allow the workbook to keep its existing sheet names (mirrors the author's
# commit message: uploaded file is read by pandas, sheet names must stay put)
$wb = $excel.ActiveWorkbook

# Rename Sheet1 -> DataSheet (same sheetId/r:id, just renamed)
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "DataSheet"

# Add a second sheet, "Instructions", immediately after DataSheet
$wsInstr = $wb.Worksheets.Add($null, $ws1)
$wsInstr.Name = "Instructions"
$wsInstr.Range("A1").Value = "This is just a test file, leave all sheet names as they are"

# Leave DataSheet as the active/selected sheet, with C5 selected
$ws1.Activate() | Out-Null
$ws1.Range("C5").Select() | Out-Null
